$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Reference cell with the default (unstyled) format, used to restore
# style on numeric-looking text cells after forcing a text NumberFormat
# (mirrors the plain, un-styled text cells in the source workbook).
$refStyle = $ws.Range("B2")

$r = $ws.Range('D2')
$r.NumberFormat = "@"
$r.Value = '22.480.96'
$r.Style = $refStyle.Style
$ws.Range('E2').Value = '  +0.41%  '
$r = $ws.Range('D3')
$r.NumberFormat = "@"
$r.Value = '1.575.84'
$r.Style = $refStyle.Style
$ws.Range('E3').Value = '  +1.03%  '
$r = $ws.Range('D4')
$r.NumberFormat = "@"
$r.Value = '1.000'
$r.Style = $refStyle.Style
$ws.Range('E5').Value = '  -0.21%  '
$r = $ws.Range('D6')
$r.NumberFormat = "@"
$r.Value = '288.28'
$r.Style = $refStyle.Style
$ws.Range('E6').Value = '  +0.75%  '
$r = $ws.Range('D7')
$r.NumberFormat = "@"
$r.Value = '0.3682'
$r.Style = $refStyle.Style
$ws.Range('E7').Value = '  +1.00%  '
$r = $ws.Range('D8')
$r.NumberFormat = "@"
$r.Value = '47.87'
$r.Style = $refStyle.Style
$ws.Range('E8').Value = '  -1.64%  '
$r = $ws.Range('D9')
$r.NumberFormat = "@"
$r.Value = '0.3330'
$r.Style = $refStyle.Style
$ws.Range('E9').Value = '  -0.29%  '
$ws.Range('E10').Value = '  +2.54%  '
$r = $ws.Range('D11')
$r.NumberFormat = "@"
$r.Value = '0.07561'
$r.Style = $refStyle.Style
$ws.Range('E11').Value = '  +2.44%  '
$ws.Range('E12').Value = '  -0.14%  '
$r = $ws.Range('D13')
$r.NumberFormat = "@"
$r.Value = '20.84'
$r.Style = $refStyle.Style
$ws.Range('E13').Value = '  +0.66%  '
$r = $ws.Range('D14')
$r.NumberFormat = "@"
$r.Value = '5.950'
$r.Style = $refStyle.Style
$ws.Range('E14').Value = '  +0.92%  '
$r = $ws.Range('D15')
$r.NumberFormat = "@"
$r.Value = '6.955'
$r.Style = $refStyle.Style
$ws.Range('E15').Value = '  +1.50%  '
$r = $ws.Range('D16')
$r.NumberFormat = "@"
$r.Value = '1.569.67'
$r.Style = $refStyle.Style
$ws.Range('E16').Value = '  +0.55%  '
$ws.Range('E17').Value = '  +2.34%  '
$r = $ws.Range('D18')
$r.NumberFormat = "@"
$r.Value = '88.01'
$r.Style = $refStyle.Style
$ws.Range('E18').Value = '  -0.87%  '
$r = $ws.Range('D19')
$r.NumberFormat = "@"
$r.Value = '0.06740'
$r.Style = $refStyle.Style
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$r = $ws.Range('D20')
$r.NumberFormat = "@"
$r.Value = '6.399'
$r.Style = $refStyle.Style
$ws.Range('E20').Value = '  +2.02%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$r = $ws.Range('D21')
$r.NumberFormat = "@"
$r.Value = '0.9997'
$r.Style = $refStyle.Style
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('E22').Value = '  +3.28%  '
$r = $ws.Range('D23')
$r.NumberFormat = "@"
$r.Value = '12.03'
$r.Style = $refStyle.Style
$ws.Range('E23').Value = '  +0.96%  '
$r = $ws.Range('D24')
$r.NumberFormat = "@"
$r.Value = '22.463.67'
$r.Style = $refStyle.Style
$ws.Range('E24').Value = '  +0.34%  '
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('E26').Value = '  +4.76%  '
$r = $ws.Range('D27')
$r.NumberFormat = "@"
$r.Value = '150.58'
$r.Style = $refStyle.Style
$ws.Range('E27').Value = '  +0.72%  '
$r = $ws.Range('D28')
$r.NumberFormat = "@"
$r.Value = '19.68'
$r.Style = $refStyle.Style
$ws.Range('E28').Value = '  +1.30%  '
$r = $ws.Range('D29')
$r.NumberFormat = "@"
$r.Value = '4.991'
$r.Style = $refStyle.Style
$ws.Range('E29').Value = '  -0.10%  '
$r = $ws.Range('D30')
$r.NumberFormat = "@"
$r.Value = '125.61'
$r.Style = $refStyle.Style
$ws.Range('E30').Value = '  +2.27%  '
$r = $ws.Range('D31')
$r.NumberFormat = "@"
$r.Value = '1.748.39'
$r.Style = $refStyle.Style
$ws.Range('E31').Value = '  +0.72%  '
$r = $ws.Range('D32')
$r.NumberFormat = "@"
$r.Value = '1.090'
$r.Style = $refStyle.Style
$r = $ws.Range('D33')
$r.NumberFormat = "@"
$r.Value = '6.112'
$r.Style = $refStyle.Style
$ws.Range('E33').Value = '  +0.41%  '
$r = $ws.Range('D34')
$r.NumberFormat = "@"
$r.Value = '1.999'
$r.Style = $refStyle.Style
$r = $ws.Range('D35')
$r.NumberFormat = "@"
$r.Value = '9.868'
$r.Style = $refStyle.Style
$ws.Range('E35').Value = '  +3.59%  '
$r = $ws.Range('D36')
$r.NumberFormat = "@"
$r.Value = '0.08372'
$r.Style = $refStyle.Style
$ws.Range('E36').Value = '  +1.71%  '
$r = $ws.Range('D37')
$r.NumberFormat = "@"
$r.Value = '0.02466'
$r.Style = $refStyle.Style
$ws.Range('E37').Value = '  +3.76%  '
$r = $ws.Range('D38')
$r.NumberFormat = "@"
$r.Value = '0.2239'
$r.Style = $refStyle.Style
$ws.Range('E38').Value = '  +1.36%  '
$r = $ws.Range('D39')
$r.NumberFormat = "@"
$r.Value = '0.06400'
$r.Style = $refStyle.Style
$ws.Range('E39').Value = '  +0.48%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$r = $ws.Range('D40')
$r.NumberFormat = "@"
$r.Value = '1.297'
$r.Style = $refStyle.Style
$ws.Range('E40').Value = '  -0.44%  '
$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$r = $ws.Range('D41')
$r.NumberFormat = "@"
$r.Value = '5.372'
$r.Style = $refStyle.Style
$ws.Range('E41').Value = '  +1.26%  '
$r = $ws.Range('D42')
$r.NumberFormat = "@"
$r.Value = '11.49'
$r.Style = $refStyle.Style
$ws.Range('E42').Value = '  +3.49%  '
$r = $ws.Range('D43')
$r.NumberFormat = "@"
$r.Value = '0.6276'
$r.Style = $refStyle.Style
$ws.Range('E43').Value = '  +3.93%  '
$r = $ws.Range('D44')
$r.NumberFormat = "@"
$r.Value = '14.16'
$r.Style = $refStyle.Style
$ws.Range('E44').Value = '  +4.09%  '
$r = $ws.Range('D45')
$r.NumberFormat = "@"
$r.Value = '0.9999'
$r.Style = $refStyle.Style
$ws.Range('E45').Value = '  -0.21%  '
$r = $ws.Range('D46')
$r.NumberFormat = "@"
$r.Value = '0.6127'
$r.Style = $refStyle.Style
$ws.Range('E46').Value = '  +7.10%  '
$r = $ws.Range('D47')
$r.NumberFormat = "@"
$r.Value = '3.779'
$r.Style = $refStyle.Style
$ws.Range('E47').Value = '  +0.51%  '
$ws.Range('E48').Value = '  +2.78%  '
$r = $ws.Range('D49')
$r.NumberFormat = "@"
$r.Value = '125.33'
$r.Style = $refStyle.Style
$ws.Range('E49').Value = '  +0.70%  '
$r = $ws.Range('D50')
$r.NumberFormat = "@"
$r.Value = '1.215'
$r.Style = $refStyle.Style
$ws.Range('E50').Value = '  +0.84%  '
$r = $ws.Range('D51')
$r.NumberFormat = "@"
$r.Value = '0.07225'
$r.Style = $refStyle.Style
$ws.Range('E51').Value = '  +0.12%  '
